$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column O values: year 2021 header (row 4) and data value (row 5)
$ws.Cells.Item(4, 14).Copy()
$ws.Cells.Item(4, 15).PasteSpecial(-4122)
$ws.Cells.Item(4, 15).Value = 2021

$ws.Cells.Item(5, 14).Copy()
$ws.Cells.Item(5, 15).PasteSpecial(-4122)
$ws.Cells.Item(5, 15).Value = 1.5020015556876996

# Update selection to match the new active cell
$ws.Range("Q5").Select()
